$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Insert 6 blank rows starting at row 5 - this pushes the existing
# table (old rows 5-11) down to rows 11-17, opening up space for the
# new "Source Type" line at row 9.
$ws.Range("A5:A10").EntireRow.Insert()

# Re-assert the pre-existing formatting on cells that were already in
# the workbook (large "name" title + bold "title" header/labels),
# since round-tripping the workbook otherwise drops their styling.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true
$ws.Range("B11:D11").Font.Bold = $true
$ws.Range("A12:A16").Font.Bold = $true

# Re-assert the italic "source" formatting on the existing byline.
$ws.Range("A17").Font.Italic = $true

# New bold+underlined "title_" style line above the data table.
$c = $ws.Range("A9")
$c.Value = "Source Type: SME Associations (Most Widely Used)"
$c.Font.Bold = $true
$c.Font.Underline = $true

# New data point: Employment (% of total) for MSMEs, which previously
# had no value. Force text storage (matches the other data cells,
# which are all shared-string text, not numbers) and then drop back to
# the default "Normal" style so it doesn't pick up a stray number
# format style.
$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "94.8"
$d.Style = "Normal"

# New source-detail rows appended below the existing source line
# (row 17), after a gap, mirroring the "title"/"source" style pairing
# already used elsewhere ("AFDB" bold like the other row labels, the
# citation italicized like the "Source: AFDB, 2006" byline).
$a = $ws.Range("A23")
$a.Value = "AFDB"
$a.Font.Bold = $true

$b = $ws.Range("A24")
$b.Value = 'African Development Bank, "Egypt Private Sector Country Profile", 2009, p. 115, 113. Available at http://www.afdb.org/fileadmin/uploads/afdb/Documents/Project-and-Operations/Brochure%20Egypt%20Anglais.pdf'
$b.Font.Italic = $true
